$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the now-superfluous data rows (old rows 5-7 corresponded to a
# sending cluster of FAPs that is being folded into rows 2-4 with
# recalculated TPM-derived values); final table only has rows 1 (header)
# .. 4 (data).
$ws.Rows("5:7").Delete()

# Row 2
$ws.Range("A2").Value = "FAPs"
$ws.Range("B2").Value = "Lif"
$ws.Range("C2").Value = "Lifr"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 2.482764666666667
$ws.Range("H2").Value = 7.448294000000001
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 1
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 34.53319033333333
$ws.Range("N2").Value = 103.599571
$ws.Range("O2").Value = 0.2461870921144496
$ws.Range("P2").Value = 0.2461870921144496
$ws.Range("Q2").Value = 85.73778478687488
$ws.Range("R2").Value = 771.640063081874
$ws.Range("S2").Value = 0.2461870921144496
$ws.Range("T2").Value = 0.2461870921144496

# Row 3
$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Lif"
$ws.Range("C3").Value = "Lifr"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 2.482764666666667
$ws.Range("H3").Value = 7.448294000000001
$ws.Range("I3").Value = 1
$ws.Range("J3").Value = 1
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 80.77474466666666
$ws.Range("N3").Value = 242.324234
$ws.Range("O3").Value = 0.575843103803214
$ws.Range("P3").Value = 0.575843103803214
$ws.Range("Q3").Value = 200.5446820174218
$ws.Range("R3").Value = 1804.902138156796
$ws.Range("S3").Value = 0.575843103803214
$ws.Range("T3").Value = 0.575843103803214

# Row 4
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Lif"
$ws.Range("C4").Value = "Lifr"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 2.482764666666667
$ws.Range("H4").Value = 7.448294000000001
$ws.Range("I4").Value = 1
$ws.Range("J4").Value = 1
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 24.96420533333334
$ws.Range("N4").Value = 74.892616
$ws.Range("O4").Value = 0.1779698040823365
$ws.Range("P4").Value = 0.1779698040823364
$ws.Range("Q4").Value = 61.98024693301156
$ws.Range("R4").Value = 557.822222397104
$ws.Range("S4").Value = 0.1779698040823365
$ws.Range("T4").Value = 0.1779698040823364
